$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a2 = New-Object 'object[,]' 1,4
$a2[0,0] = 144.783305
$a2[0,1] = 434.349915
$a2[0,2] = 0.2430046335191003
$a2[0,3] = 0.251012682214973
$ws.Range("G2:J2").Value = $a2

$b2 = New-Object 'object[,]' 1,8
$b2[0,0] = 8.021311666666668
$b2[0,1] = 24.063935
$b2[0,2] = 0.2318437811880539
$b2[0,3] = 0.2575223477274881
$b2[0,4] = 1161.352013535058
$b2[0,5] = 10452.16812181552
$b2[0,6] = 0.05633911308128552
$b2[0,7] = 0.06464137523337374
$ws.Range("M2:T2").Value = $b2

$a3 = New-Object 'object[,]' 1,4
$a3[0,0] = 144.783305
$a3[0,1] = 434.349915
$a3[0,2] = 0.2430046335191003
$a3[0,3] = 0.251012682214973
$ws.Range("G3:J3").Value = $a3

$b3 = New-Object 'object[,]' 1,8
$b3[0,0] = 11.70102866666667
$b3[0,1] = 35.103086
$b3[0,2] = 0.3382003894878138
$b3[0,3] = 0.3756588072233373
$b3[0,4] = 1694.113602259743
$b3[0,5] = 15247.02242033769
$b3[0,6] = 0.08218426170350319
$b3[0,7] = 0.09429512479880735
$ws.Range("M3:T3").Value = $b3

$a4 = New-Object 'object[,]' 1,4
$a4[0,0] = 144.783305
$a4[0,1] = 434.349915
$a4[0,2] = 0.2430046335191003
$a4[0,3] = 0.251012682214973
$ws.Range("G4:J4").Value = $a4

$b4 = New-Object 'object[,]' 1,8
$b4[0,0] = 2.475956
$b4[0,1] = 7.427868
$b4[0,2] = 0.07156373233578578
$b4[0,3] = 0.07948999222154987
$b4[0,4] = 358.47709271458
$b4[0,5] = 3226.29383443122
$b4[0,6] = 0.01739031854951661
$b4[0,7] = 0.01995299615677857
$ws.Range("M4:T4").Value = $b4

$a5 = New-Object 'object[,]' 1,4
$a5[0,0] = 144.783305
$a5[0,1] = 434.349915
$a5[0,2] = 0.2430046335191003
$a5[0,3] = 0.251012682214973
$ws.Range("G5:J5").Value = $a5

$b5 = New-Object 'object[,]' 1,8
$b5[0,0] = 10.349683
$b5[0,1] = 20.699366
$b5[0,2] = 0.2991418038011306
$b5[0,3] = 0.2215161123395049
$b5[0,4] = 1498.461310442315
$b5[0,5] = 8990.767862653887
$b5[0,6] = 0.07269284440293637
$b5[0,7] = 0.0556033535121724
$ws.Range("M5:T5").Value = $b5

$a6 = New-Object 'object[,]' 1,4
$a6[0,0] = 144.783305
$a6[0,1] = 434.349915
$a6[0,2] = 0.2430046335191003
$a6[0,3] = 0.251012682214973
$ws.Range("G6:J6").Value = $a6

$b6 = New-Object 'object[,]' 1,8
$b6[0,0] = 2.049936666666667
$b6[0,1] = 6.14981
$b6[0,2] = 0.05925029318721588
$b6[0,3] = 0.06581274048811983
$b6[0,4] = 296.7966056406834
$b6[0,5] = 2671.16945076615
$b6[0,6] = 0.01439809578185864
$b6[0,7] = 0.01651983251384091
$ws.Range("M6:T6").Value = $b6

$a7 = New-Object 'object[,]' 1,4
$a7[0,0] = 82.24887099999999
$a7[0,1] = 246.746613
$a7[0,2] = 0.1380466950572427
$a7[0,3] = 0.1425959278859072
$ws.Range("G7:J7").Value = $a7

$b7 = New-Object 'object[,]' 1,8
$b7[0,0] = 8.021311666666668
$b7[0,1] = 24.063935
$b7[0,2] = 0.2318437811880539
$b7[0,3] = 0.2575223477274881
$b7[0,4] = 659.7438285224617
$b7[0,5] = 5937.694456702156
$b7[0,6] = 0.03200526776258537
$b7[0,7] = 0.03672163812555841
$ws.Range("M7:T7").Value = $b7

$a8 = New-Object 'object[,]' 1,4
$a8[0,0] = 82.24887099999999
$a8[0,1] = 246.746613
$a8[0,2] = 0.1380466950572427
$a8[0,3] = 0.1425959278859072
$ws.Range("G8:J8").Value = $a8

$b8 = New-Object 'object[,]' 1,8
$b8[0,0] = 11.70102866666667
$b8[0,1] = 35.103086
$b8[0,2] = 0.3382003894878138
$b8[0,3] = 0.3756588072233373
$b8[0,4] = 962.3963973719685
$b8[0,5] = 8661.567576347717
$b8[0,6] = 0.04668744603586494
$b8[0,7] = 0.05356741618452492
$ws.Range("M8:T8").Value = $b8

$a9 = New-Object 'object[,]' 1,4
$a9[0,0] = 82.24887099999999
$a9[0,1] = 246.746613
$a9[0,2] = 0.1380466950572427
$a9[0,3] = 0.1425959278859072
$ws.Range("G9:J9").Value = $a9

$b9 = New-Object 'object[,]' 1,8
$b9[0,0] = 2.475956
$b9[0,1] = 7.427868
$b9[0,2] = 0.07156373233578578
$b9[0,3] = 0.07948999222154987
$b9[0,4] = 203.644585645676
$b9[0,5] = 1832.801270811084
$b9[0,6] = 0.009879136734916357
$b9[0,7] = 0.01133494919847545
$ws.Range("M9:T9").Value = $b9

$a10 = New-Object 'object[,]' 1,4
$a10[0,0] = 82.24887099999999
$a10[0,1] = 246.746613
$a10[0,2] = 0.1380466950572427
$a10[0,3] = 0.1425959278859072
$ws.Range("G10:J10").Value = $a10

$b10 = New-Object 'object[,]' 1,8
$b10[0,0] = 10.349683
$b10[0,1] = 20.699366
$b10[0,2] = 0.2991418038011306
$b10[0,3] = 0.2215161123395049
$b10[0,4] = 851.2497419578929
$b10[0,5] = 5107.498451747358
$b10[0,6] = 0.0412955373682082
$b10[0,7] = 0.03158729558073056
$ws.Range("M10:T10").Value = $b10

$a11 = New-Object 'object[,]' 1,4
$a11[0,0] = 82.24887099999999
$a11[0,1] = 246.746613
$a11[0,2] = 0.1380466950572427
$a11[0,3] = 0.1425959278859072
$ws.Range("G11:J11").Value = $a11

$b11 = New-Object 'object[,]' 1,8
$b11[0,0] = 2.049936666666667
$b11[0,1] = 6.14981
$b11[0,2] = 0.05925029318721588
$b11[0,3] = 0.06581274048811983
$b11[0,4] = 168.6049764548367
$b11[0,5] = 1517.44478809353
$b11[0,6] = 0.008179307155667814
$b11[0,7] = 0.009384628796617859
$ws.Range("M11:T11").Value = $b11

$a12 = New-Object 'object[,]' 1,4
$a12[0,0] = 163.8590903333333
$a12[0,1] = 491.577271
$a12[0,2] = 0.2750214756820535
$a12[0,3] = 0.284084617144743
$ws.Range("G12:J12").Value = $a12

$b12 = New-Object 'object[,]' 1,8
$b12[0,0] = 8.021311666666668
$b12[0,1] = 24.063935
$b12[0,2] = 0.2318437811880539
$b12[0,3] = 0.2575223477274881
$b12[0,4] = 1314.364832980154
$b12[0,5] = 11829.28349682138
$b12[0,6] = 0.06376201883004569
$b12[0,7] = 0.07315813756037884
$ws.Range("M12:T12").Value = $b12

$a13 = New-Object 'object[,]' 1,4
$a13[0,0] = 163.8590903333333
$a13[0,1] = 491.577271
$a13[0,2] = 0.2750214756820535
$a13[0,3] = 0.284084617144743
$ws.Range("G13:J13").Value = $a13

$b13 = New-Object 'object[,]' 1,8
$b13[0,0] = 11.70102866666667
$b13[0,1] = 35.103086
$b13[0,2] = 0.3382003894878138
$b13[0,3] = 0.3756588072233373
$b13[0,4] = 1917.319913284256
$b13[0,5] = 17255.8792195583
$b13[0,6] = 0.09301237019318381
$b13[0,7] = 0.1067188884270926
$ws.Range("M13:T13").Value = $b13

$a14 = New-Object 'object[,]' 1,4
$a14[0,0] = 163.8590903333333
$a14[0,1] = 491.577271
$a14[0,2] = 0.2750214756820535
$a14[0,3] = 0.284084617144743
$ws.Range("G14:J14").Value = $a14

$b14 = New-Object 'object[,]' 1,8
$b14[0,0] = 2.475956
$b14[0,1] = 7.427868
$b14[0,2] = 0.07156373233578578
$b14[0,3] = 0.07948999222154987
$b14[0,4] = 405.7078978653587
$b14[0,5] = 3651.371080788228
$b14[0,6] = 0.01968156327230329
$b14[0,7] = 0.02258188400709759
$ws.Range("M14:T14").Value = $b14

$a15 = New-Object 'object[,]' 1,4
$a15[0,0] = 163.8590903333333
$a15[0,1] = 491.577271
$a15[0,2] = 0.2750214756820535
$a15[0,3] = 0.284084617144743
$ws.Range("G15:J15").Value = $a15

$b15 = New-Object 'object[,]' 1,8
$b15[0,0] = 10.349683
$b15[0,1] = 20.699366
$b15[0,2] = 0.2991418038011306
$b15[0,3] = 0.2215161123395049
$b15[0,4] = 1695.889641618364
$b15[0,5] = 10175.33784971018
$b15[0,6] = 0.08227042031957826
$b15[0,7] = 0.06292931996536012
$ws.Range("M15:T15").Value = $b15

$a16 = New-Object 'object[,]' 1,4
$a16[0,0] = 163.8590903333333
$a16[0,1] = 491.577271
$a16[0,2] = 0.2750214756820535
$a16[0,3] = 0.284084617144743
$ws.Range("G16:J16").Value = $a16

$b16 = New-Object 'object[,]' 1,8
$b16[0,0] = 2.049936666666667
$b16[0,1] = 6.14981
$b16[0,2] = 0.05925029318721588
$b16[0,3] = 0.06581274048811983
$b16[0,4] = 335.9007574409456
$b16[0,5] = 3023.10681696851
$b16[0,6] = 0.01629510306694243
$b16[0,7] = 0.01869638718481385
$ws.Range("M16:T16").Value = $b16

$a17 = New-Object 'object[,]' 1,4
$a17[0,0] = 57.0238095
$a17[0,1] = 114.047619
$a17[0,2] = 0.09570889357312636
$a17[0,3] = 0.06590860906562239
$ws.Range("G17:J17").Value = $a17

$b17 = New-Object 'object[,]' 1,8
$b17[0,0] = 8.021311666666668
$b17[0,1] = 24.063935
$b17[0,2] = 0.2318437811880539
$b17[0,3] = 0.2575223477274881
$b17[0,4] = 457.4057484201276
$b17[0,5] = 2744.434490520765
$b17[0,6] = 0.02218951177931864
$b17[0,7] = 0.01697293974203229
$ws.Range("M17:T17").Value = $b17

$a18 = New-Object 'object[,]' 1,4
$a18[0,0] = 57.0238095
$a18[0,1] = 114.047619
$a18[0,2] = 0.09570889357312636
$a18[0,3] = 0.06590860906562239
$ws.Range("G18:J18").Value = $a18

$b18 = New-Object 'object[,]' 1,8
$b18[0,0] = 11.70102866666667
$b18[0,1] = 35.103086
$b18[0,2] = 0.3382003894878138
$b18[0,3] = 0.3756588072233373
$b18[0,4] = 667.237229642039
$b18[0,5] = 4003.423377852233
$b18[0,6] = 0.03236878508387905
$b18[0,7] = 0.02475914946734094
$ws.Range("M18:T18").Value = $b18

$a19 = New-Object 'object[,]' 1,4
$a19[0,0] = 57.0238095
$a19[0,1] = 114.047619
$a19[0,2] = 0.09570889357312636
$a19[0,3] = 0.06590860906562239
$ws.Range("G19:J19").Value = $a19

$b19 = New-Object 'object[,]' 1,8
$b19[0,0] = 2.475956
$b19[0,1] = 7.427868
$b19[0,2] = 0.07156373233578578
$b19[0,3] = 0.07948999222154987
$b19[0,4] = 141.188443274382
$b19[0,5] = 847.130659646292
$b19[0,6] = 0.006849285641821422
$b19[0,7] = 0.005239074821959495
$ws.Range("M19:T19").Value = $b19

$a20 = New-Object 'object[,]' 1,4
$a20[0,0] = 57.0238095
$a20[0,1] = 114.047619
$a20[0,2] = 0.09570889357312636
$a20[0,3] = 0.06590860906562239
$ws.Range("G20:J20").Value = $a20

$b20 = New-Object 'object[,]' 1,8
$b20[0,0] = 10.349683
$b20[0,1] = 20.699366
$b20[0,2] = 0.2991418038011306
$b20[0,3] = 0.2215161123395049
$b20[0,4] = 590.1783517773885
$b20[0,5] = 2360.713407109554
$b20[0,6] = 0.02863053106327546
$b20[0,7] = 0.01459981884992092
$ws.Range("M20:T20").Value = $b20

$a21 = New-Object 'object[,]' 1,4
$a21[0,0] = 57.0238095
$a21[0,1] = 114.047619
$a21[0,2] = 0.09570889357312636
$a21[0,3] = 0.06590860906562239
$ws.Range("G21:J21").Value = $a21

$b21 = New-Object 'object[,]' 1,8
$b21[0,0] = 2.049936666666667
$b21[0,1] = 6.14981
$b21[0,2] = 0.05925029318721588
$b21[0,3] = 0.06581274048811983
$b21[0,4] = 116.895197967065
$b21[0,5] = 701.3711878023901
$b21[0,6] = 0.005670780004831778
$b21[0,7] = 0.004337626184368748
$ws.Range("M21:T21").Value = $b21

$a22 = New-Object 'object[,]' 1,4
$a22[0,0] = 147.8896333333333
$a22[0,1] = 443.6689
$a22[0,2] = 0.2482183021684772
$a22[0,3] = 0.2563981636887546
$ws.Range("G22:J22").Value = $a22

$b22 = New-Object 'object[,]' 1,8
$b22[0,0] = 8.021311666666668
$b22[0,1] = 24.063935
$b22[0,2] = 0.2318437811880539
$b22[0,3] = 0.2575223477274881
$b22[0,4] = 1186.268841235722
$b22[0,5] = 10676.4195711215
$b22[0,6] = 0.05754786973481867
$b22[0,7] = 0.06602825706614487
$ws.Range("M22:T22").Value = $b22

$a23 = New-Object 'object[,]' 1,4
$a23[0,0] = 147.8896333333333
$a23[0,1] = 443.6689
$a23[0,2] = 0.2482183021684772
$a23[0,3] = 0.2563981636887546
$ws.Range("G23:J23").Value = $a23

$b23 = New-Object 'object[,]' 1,8
$b23[0,0] = 11.70102866666667
$b23[0,1] = 35.103086
$b23[0,2] = 0.3382003894878138
$b23[0,3] = 0.3756588072233373
$b23[0,4] = 1730.460839136156
$b23[0,5] = 15574.1475522254
$b23[0,6] = 0.08394752647138286
$b23[0,7] = 0.09631822834557152
$ws.Range("M23:T23").Value = $b23

$a24 = New-Object 'object[,]' 1,4
$a24[0,0] = 147.8896333333333
$a24[0,1] = 443.6689
$a24[0,2] = 0.2482183021684772
$a24[0,3] = 0.2563981636887546
$ws.Range("G24:J24").Value = $a24

$b24 = New-Object 'object[,]' 1,8
$b24[0,0] = 2.475956
$b24[0,1] = 7.427868
$b24[0,2] = 0.07156373233578578
$b24[0,3] = 0.07948999222154987
$b24[0,4] = 366.1682249894667
$b24[0,5] = 3295.5140249052
$b24[0,6] = 0.0177634281372281
$b24[0,7] = 0.02038108803723877
$ws.Range("M24:T24").Value = $b24

$a25 = New-Object 'object[,]' 1,4
$a25[0,0] = 147.8896333333333
$a25[0,1] = 443.6689
$a25[0,2] = 0.2482183021684772
$a25[0,3] = 0.2563981636887546
$ws.Range("G25:J25").Value = $a25

$b25 = New-Object 'object[,]' 1,8
$b25[0,0] = 10.349683
$b25[0,1] = 20.699366
$b25[0,2] = 0.2991418038011306
$b25[0,3] = 0.2215161123395049
$b25[0,4] = 1530.610823986233
$b25[0,5] = 9183.6649439174
$b25[0,6] = 0.07425247064713238
$b25[0,7] = 0.05679632443132093
$ws.Range("M25:T25").Value = $b25

$a26 = New-Object 'object[,]' 1,4
$a26[0,0] = 147.8896333333333
$a26[0,1] = 443.6689
$a26[0,2] = 0.2482183021684772
$a26[0,3] = 0.2563981636887546
$ws.Range("G26:J26").Value = $a26

$b26 = New-Object 'object[,]' 1,8
$b26[0,0] = 2.049936666666667
$b26[0,1] = 6.14981
$b26[0,2] = 0.05925029318721588
$b26[0,3] = 0.06581274048811983
$b26[0,4] = 303.164381989889
$b26[0,5] = 2728.479437909
$b26[0,6] = 0.01470700717791522
$b26[0,7] = 0.01687426580847847
$ws.Range("M26:T26").Value = $b26

